$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Slit1-Robo2 table ("Natmi following Dr Hou advice"):
# expands from 4 data rows (FAPs/M1/M2/sCs -> ECs only) to 12 data rows,
# covering every sending/target cluster combination, with refreshed
# expression/specificity statistics.
$data = @(
    @("FAPs", "Slit1", "Robo2", "ECs", 3, 1, 0.340305, 1.020915, 0.2413792532744959, 0.2522899505114672, 3, 1, 1.252524333333333, 3.757573, 0.9099863583404593, 0.9099863583404593, 0.426240293255, 3.836162639295, 0.2196518276661979, 0.2295804133118247),
    @("FAPs", "Slit1", "Robo2", "FAPs", 3, 1, 0.340305, 1.020915, 0.2413792532744959, 0.2522899505114672, 2, 0.6666666666666666, 0.108201, 0.324603, 0.0786103960924746, 0.07861039609247461, 0.036821341305, 0.331392071745, 0.01897491870841387, 0.01983261293985725),
    @("FAPs", "Slit1", "Robo2", "M2", 3, 1, 0.340305, 1.020915, 0.2413792532744959, 0.2522899505114672, 1, 0.3333333333333333, 0.01569566666666667, 0.047087, 0.01140324556706608, 0.01140324556706608, 0.005341313845000001, 0.048071824605, 0.002752506899884117, 0.00287692425978521),
    @("M1", "Slit1", "Robo2", "ECs", 3, 1, 0.472982, 1.418946, 0.3354874068035369, 0.3506519309819567, 3, 1, 1.252524333333333, 3.757573, 0.9099863583404593, 0.9099863583404593, 0.5924214642286667, 5.331793178058, 0.3052889635862348, 0.3190884737193208),
    @("M1", "Slit1", "Robo2", "FAPs", 3, 1, 0.472982, 1.418946, 0.3354874068035369, 0.3506519309819567, 2, 0.6666666666666666, 0.108201, 0.324603, 0.0786103960924746, 0.07861039609247461, 0.051177125382, 0.460594128438, 0.02637279793286319, 0.02756488718508269),
    @("M1", "Slit1", "Robo2", "M2", 3, 1, 0.472982, 1.418946, 0.3354874068035369, 0.3506519309819567, 1, 0.3333333333333333, 0.01569566666666667, 0.047087, 0.01140324556706608, 0.01140324556706608, 0.007423767811333333, 0.066813910302, 0.003825645284438928, 0.00399857007755316),
    @("M2", "Slit1", "Robo2", "ECs", 3, 1, 0.4136363333333333, 1.240909, 0.2933933655608953, 0.3066551771687498, 3, 1, 1.252524333333333, 3.757573, 0.9099863583404593, 0.9099863583404593, 0.5180895726507778, 4.662806153857, 0.2669839602880102, 0.279052027938039),
    @("M2", "Slit1", "Robo2", "FAPs", 3, 1, 0.4136363333333333, 1.240909, 0.2933933655608953, 0.3066551771687498, 2, 0.6666666666666666, 0.108201, 0.324603, 0.0786103960924746, 0.07861039609247461, 0.04475586490299999, 0.402802784127, 0.02306376867764617, 0.0241062849410434),
    @("M2", "Slit1", "Robo2", "M2", 3, 1, 0.4136363333333333, 1.240909, 0.2933933655608953, 0.3066551771687498, 1, 0.3333333333333333, 0.01569566666666667, 0.047087, 0.01140324556706608, 0.01140324556706608, 0.006492298009222222, 0.058430682083, 0.003345636595238878, 0.003496864289667411),
    @("sCs", "Slit1", "Robo2", "ECs", 2, 1, 0.182912, 0.365824, 0.129739974361072, 0.09040294133782634, 3, 1, 1.252524333333333, 3.757573, 0.9099863583404593, 0.9099863583404593, 0.2291017308586667, 1.374610385152, 0.1180616068000164, 0.08226544337127475),
    @("sCs", "Slit1", "Robo2", "FAPs", 2, 1, 0.182912, 0.365824, 0.129739974361072, 0.09040294133782634, 2, 0.6666666666666666, 0.108201, 0.324603, 0.0786103960924746, 0.07861039609247461, 0.019791261312, 0.118747567872, 0.01019891077355137, 0.007106611026491275),
    @("sCs", "Slit1", "Robo2", "M2", 2, 1, 0.182912, 0.365824, 0.129739974361072, 0.09040294133782634, 1, 0.3333333333333333, 0.01569566666666667, 0.047087, 0.01140324556706608, 0.01140324556706608, 0.002870925781333334, 0.017225554688, 0.001479456787504161, 0.001030886940060303)
)

$numRows = $data.Count
$numCols = $data[0].Count

# Build a true 2D array so Excel COM accepts it as a single Range.Value2 write.
$arr = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $numRows - 1
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $numCols))
$targetRange.Value2 = $arr
